# Applies the "Add files via upload" revision to Linkovi.docx:
#  1. Collapses several spell-checked (proofErr-wrapped) runs back into
#     single plain runs by doing an exact Find & Replace of the full
#     paragraph text (Word's Find/Replace naturally re-homogenizes the
#     run/proofErr structure when it rewrites the matched range).
#  2. Inserts a new "Map Test 2.8" / "New Locations" / hyperlink block
#     (and relocates the lone "_GoBack" bookmark paragraph into it),
#     consuming some of the run of blank paragraphs that used to sit
#     between the "Map Test 2.7" block and the "Mape:" heading.
#  3. Cleans up the remaining proofErr-wrapped names further down
#     (RuinsArena, SnowMountain, ForestCamp, SnowCampArena,
#     SnowClearingArena, SnowCorneredArena, SnowForestArena).

$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found -> $oldText"
    }
}

# --- 1. Merge split/spell-checked runs back into single runs ---------------

Replace-Exact "First combat test: " "First combat test: "
Replace-Exact "Equipment test:" "Equipment test:"
Replace-Exact "First map test (obsolete, algorithm changed)" "First map test (obsolete, algorithm changed)"
Replace-Exact "Map Test 2.0, " "Map Test 2.0, "
Replace-Exact "improved algorithm, faster generation, objects inside zones" "improved algorithm, faster generation, objects inside zones"
Replace-Exact "Map Tesp 2.5, " "Map Tesp 2.5, "
Replace-Exact "Allows movement between world and locations" "Allows movement between world and locations"
Replace-Exact "Map Test 2.7" "Map Test 2.7"
Replace-Exact "Added transition between world and location and vice versa, movement between locations and world seperated" "Added transition between world and location and vice versa, movement between locations and world seperated"

# --- 2. Insert the new "Map Test 2.8" / "New Locations" / link block -------
# Right after the "Map Test 2.7" hyperlink paragraph there are 10 empty
# paragraphs before the "Mape:" heading. Keep the first blank, turn the
# next three into content, keep one blank, turn the next into the
# relocated _GoBack bookmark, then drop one extra blank paragraph so the
# overall paragraph count matches (10 blanks -> 9 paragraphs).

$mapeRange = $d.Content
$mapeRange.Find.Execute("Mape:") | Out-Null
$mapePara = $mapeRange.Paragraphs.First
$mapeIndex = $mapePara.Range.Information(3)  # wdActiveEndAdjustedPageNumber not needed; fallback below

# Locate paragraphs by walking from the known "Map Test 2.7" hyperlink
# paragraph instead of relying on fixed indices, so this keeps working
# even if earlier replacements change paragraph counts.
$hyperlinkRange = $d.Content
$hyperlinkRange.Find.Execute("https://drive.google.com/open?id=0B0dYxrDwUlTxckNUeW4zTGlfUEk") | Out-Null
$anchorPara = $hyperlinkRange.Paragraphs.First
$anchorIndex = $anchorPara.Index

$p1 = $d.Paragraphs.Item($anchorIndex + 2)
$p1.Range.Text = "Map Test 2.8"

$p2 = $d.Paragraphs.Item($anchorIndex + 3)
$p2.Range.Text = "New Locations"

$p3 = $d.Paragraphs.Item($anchorIndex + 4)
$r3 = $p3.Range
$r3.Collapse(1)
$d.Hyperlinks.Add($r3, "https://drive.google.com/open?id=0B0dYxrDwUlTxcGUwMVE4LUJBSDg", $null, $null, "https://drive.google.com/open?id=0B0dYxrDwUlTxcGUwMVE4LUJBSDg") | Out-Null

# paragraph (anchorIndex + 5) stays blank

$p5 = $d.Paragraphs.Item($anchorIndex + 6)
$d.Bookmarks.Add("_GoBack", $p5.Range) | Out-Null

# remove one extra blank paragraph so the 10 originally-blank paragraphs
# become the target 9 (1 blank + 3 content + 1 blank + 1 bookmark + 3 blank)
$d.Paragraphs.Item($anchorIndex + 7).Range.Delete()

# --- 3. Clean up remaining proofErr-wrapped names ---------------------------

Replace-Exact "RuinsArena, Ratomir" "RuinsArena, Ratomir"
Replace-Exact "SnowMountain, Ratomir" "SnowMountain, Ratomir"
Replace-Exact "ForestCamp, Ratomir" "ForestCamp, Ratomir"
Replace-Exact "SnowCampArena, Toni" "SnowCampArena, Toni"
Replace-Exact "SnowClearingArena, Toni" "SnowClearingArena, Toni"
Replace-Exact "SnowCorneredArena, Toni" "SnowCorneredArena, Toni"
Replace-Exact "SnowForestArena, Toni" "SnowForestArena, Toni"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
